$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 640, shifting existing rows 640:720 down to 642:722
$ws.Rows("640:641").Insert()

# Fill in the new row 640 (Fukumoto / Primera)
$ws.Range("A640").Value = 7
$ws.Range("B640").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C640").Value = "Ñuble"
$ws.Range("D640").Value = 45127
$ws.Range("E640").Value = 16
$ws.Range("F640").Value = "Fruta"
$ws.Range("G640").Value = 100102
$ws.Range("H640").Value = "Cítricos"
$ws.Range("I640").Value = 100102005
$ws.Range("J640").Value = "Naranja"
$ws.Range("K640").Value = "Fukumoto"
$ws.Range("L640").Value = "Primera"
$ws.Range("M640").Value = 80
$ws.Range("N640").Value = 10000
$ws.Range("O640").Value = 10000
$ws.Range("P640").Value = 10000
$ws.Range("Q640").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R640").Value = "Región de O'Higgins"
$ws.Range("S640").Value = 667
$ws.Range("T640").Value = 15

# Fill in the new row 641 (Fukumoto / Segunda)
$ws.Range("A641").Value = 7
$ws.Range("B641").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C641").Value = "Ñuble"
$ws.Range("D641").Value = 45127
$ws.Range("E641").Value = 16
$ws.Range("F641").Value = "Fruta"
$ws.Range("G641").Value = 100102
$ws.Range("H641").Value = "Cítricos"
$ws.Range("I641").Value = 100102005
$ws.Range("J641").Value = "Naranja"
$ws.Range("K641").Value = "Fukumoto"
$ws.Range("L641").Value = "Segunda"
$ws.Range("M641").Value = 60
$ws.Range("N641").Value = 8000
$ws.Range("O641").Value = 8000
$ws.Range("P641").Value = 8000
$ws.Range("Q641").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R641").Value = "Región de O'Higgins"
$ws.Range("S641").Value = 533
$ws.Range("T641").Value = 15
